$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating-point precision of the existing last row's date value
$ws.Range("A29").Value = 44342.79180781945

# Append the newly retrieved row of job-number data
$ws.Range("A30").Value = 44343.79602429144
$ws.Range("B30").Value = 74460
$ws.Range("C30").Value = 62627
$ws.Range("D30").Value = 3315
$ws.Range("E30").Value = 2055
$ws.Range("F30").Value = 1430
$ws.Range("G30").Value = 19383
$ws.Range("H30").Value = 1391
$ws.Range("I30").Value = 821
$ws.Range("J30").Value = 213
